# Generate Report for Handback
#
# For both the zh-cn and de-de localization-status worksheets:
#  - Status (col C) moves from "Ready for handoff" to "Handed back: in sync with en-US"
#  - Two new columns get populated for the already-handed-off rows:
#      F "Latest Target File"   -> same source .md file/link as column A
#      G "Latest Handback File" -> same .xlf file/link as column D
#  - Latest Handback DateTime (col H) is stamped with the handback time
#    (different per-language, since each language finishes at a different time)

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "Handed back: in sync with en-US"

function Update-HandbackSheet($SheetName, $HandbackStamp) {

    $ws = $wb.Worksheets.Item($SheetName)

    # Capture the existing hyperlink targets for A2/D2/A3/D3 before we add
    # any new ones, keyed off the cell address they're anchored to.
    $a2url = ""
    $d2url = ""
    $a3url = ""
    $d3url = ""
    foreach ($hl in $ws.Hyperlinks) {
        $addr = $hl.Range.Address()
        if ($addr -eq "`$A`$2") { $a2url = $hl.Address }
        if ($addr -eq "`$D`$2") { $d2url = $hl.Address }
        if ($addr -eq "`$A`$3") { $a3url = $hl.Address }
        if ($addr -eq "`$D`$3") { $d3url = $hl.Address }
    }

    $a2val = $ws.Range("A2").Value()
    $d2val = $ws.Range("D2").Value()
    $a3val = $ws.Range("A3").Value()
    $d3val = $ws.Range("D3").Value()

    # Status -> handed back
    if ($ws.Range("C2").Value() -eq $oldStatus) {
        $ws.Range("C2").Value = $newStatus
    }
    if ($ws.Range("C3").Value() -eq $oldStatus) {
        $ws.Range("C3").Value = $newStatus
    }

    # New "Latest Target File" (F) / "Latest Handback File" (G) columns,
    # mirroring the handoff file (A/md) and handoff target (D/xlf) links.
    $ws.Range("F2").Value = $a2val
    $ws.Hyperlinks.Add($ws.Range("F2"), $a2url, "", "", $a2val) | Out-Null
    $ws.Range("F2").Style = "HyperLink"

    $ws.Range("G2").Value = $d2val
    $ws.Hyperlinks.Add($ws.Range("G2"), $d2url, "", "", $d2val) | Out-Null
    $ws.Range("G2").Style = "HyperLink"

    $ws.Range("F3").Value = $a3val
    $ws.Hyperlinks.Add($ws.Range("F3"), $a3url, "", "", $a3val) | Out-Null
    $ws.Range("F3").Style = "HyperLink"

    $ws.Range("G3").Value = $d3val
    $ws.Hyperlinks.Add($ws.Range("G3"), $d3url, "", "", $d3val) | Out-Null
    $ws.Range("G3").Style = "HyperLink"

    # Latest Handback DateTime (H) for both rows.
    $ws.Range("H2").Value = $HandbackStamp
    $ws.Range("H3").Value = $HandbackStamp
}

Update-HandbackSheet "zh-cn" "2016-03-23 10:37:01"
Update-HandbackSheet "de-de" "2016-03-23 10:37:17"
